# Regenerate merged AHB files
# 1) Rename header row captions from *_old/_new to *_FV2210/_FV2304
# 2) Wrap the data range in a native Excel Table (ListObject)
# 3) Freeze the header row (row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseCols  = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")
$leftCols  = @("A","B","C","D","E","F","G","H","I","J")
$rightCols = @("L","M","N","O","P","Q","R","S","T","U")

for ($i = 0; $i -lt $baseCols.Length; $i++) {
    $ws.Range($leftCols[$i]  + "1").Value2 = $baseCols[$i] + "_FV2210"
    $ws.Range($rightCols[$i] + "1").Value2 = $baseCols[$i] + "_FV2304"
}
# column K ("diff") is unchanged

# --- Build the table without inheriting the header's dxf/style capture ---
# Stash the current header formatting in an unused area of the sheet so it
# survives a ClearFormats() call (ListObjects.Add bakes any pre-existing
# header formatting into a dxf + named table style, which the target file
# does not have).
$hdr = $ws.Range("A1:U1")
$stash = $ws.Range("A200:U200")
[void]$hdr.Copy()
[void]$stash.PasteSpecial(-4122)
$excel.CutCopyMode = $false

[void]$hdr.ClearFormats()

$dataRange = $ws.Range("A1:U93")
$tbl = $ws.ListObjects.Add(1, $dataRange, [System.Reflection.Missing]::Value, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# Restore the original header formatting (bold / fill / border / centered, wrapped).
[void]$stash.Copy()
[void]$hdr.PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Clean up the stash area completely.
[void]$stash.Clear()

# --- Freeze the header row ---
[void]$ws.Activate()
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
